$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of row 9 down into row 10 so the new row matches
# the existing style (number formats, etc.) used by the previous data rows.
$ws.Range("A9:E9").Copy()
$ws.Range("A10:E10").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Populate the new data row.
$ws.Range("A10").Value = "Empresarial"
$ws.Range("B10").Value = "Nova Sede ADM"
$ws.Range("C10").Value = -17.790740803423201
$ws.Range("D10").Value = -50.922744466387201
$ws.Range("E10").Value = "yellow"

# Update the selected cell as recorded in the saved view.
$ws.Range("E11").Select()
